$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Moussa Diabate", "C", "Charlotte Hornets"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Max Christie", "SG,SF", "Dallas Mavericks"),
    @("Jimmy Butler III", "SF,PF", "Golden State Warriors"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Bol Bol", "PF,C", "Phoenix Suns"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Brandon Ingram", "SG,SF,PF", "Toronto Raptors"),
    @("Devin Booker", "PG,SG", "Phoenix Suns")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
